$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 297.63043
$ws.Range("I33").Value = 262.53488
$ws.Range("J33").Value = 800.6667
$ws.Range("K33").Value = 262.53488
$ws.Range("L33").Value = 800.6667
$ws.Range("M33").Value = -33.53487999999999
$ws.Range("N33").Value = -1258.6667
$ws.Range("H39").Value = 202.09091
$ws.Range("I39").Value = 110.42857
$ws.Range("J39").Value = 362.5
$ws.Range("K39").Value = 331.28571
$ws.Range("L39").Value = 1087.5
$ws.Range("M39").Value = -35.28570999999999
$ws.Range("N39").Value = -1679.5
$ws.Range("H86").Value = 16440.857
$ws.Range("I86").Value = 2025
$ws.Range("J86").Value = 35662
$ws.Range("K86").Value = 2025
$ws.Range("L86").Value = 35662
$ws.Range("M86").Value = -902
$ws.Range("N86").Value = -37908
$ws.Range("H89").Value = 16440.857
$ws.Range("I89").Value = 2025
$ws.Range("J89").Value = 35662
$ws.Range("K89").Value = 10125
$ws.Range("L89").Value = 178310
$ws.Range("M89").Value = -4509
$ws.Range("N89").Value = -189542
$ws.Range("H98").Value = 851
$ws.Range("I98").Value = 851
$ws.Range("K98").Value = 851
$ws.Range("M98").Value = 647
$ws.Range("H106").Value = 2083.0454
$ws.Range("I106").Value = 1280.5
$ws.Range("J106").Value = 2751.8333
$ws.Range("K106").Value = 1280.5
$ws.Range("L106").Value = 2751.8333
$ws.Range("M106").Value = -649.5
$ws.Range("N106").Value = -4013.8333
$ws.Range("H113").Value = 100004184
$ws.Range("I113").Value = 142859710
$ws.Range("J113").Value = 7962.6665
$ws.Range("K113").Value = 142859710
$ws.Range("L113").Value = 7962.6665
$ws.Range("M113").Value = -142856456
$ws.Range("N113").Value = -14470.6665
$ws.Range("H122").Value = 851
$ws.Range("I122").Value = 851
$ws.Range("K122").Value = 2553
$ws.Range("M122").Value = -103
$ws.Range("H138").Value = 2767.4119
$ws.Range("I138").Value = 1794.875
$ws.Range("J138").Value = 3631.889
$ws.Range("K138").Value = 5384.625
$ws.Range("L138").Value = 10895.667
$ws.Range("M138").Value = -244.625
$ws.Range("N138").Value = -21175.667

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4501.278
$ws.Range("I32").Value = 3215.4285
$ws.Range("K32").Value = 3215.4285
$ws.Range("M32").Value = -2928.4285
$ws.Range("H61").Value = 2436.8147
$ws.Range("I61").Value = 1198.5714
$ws.Range("K61").Value = 1198.5714
$ws.Range("M61").Value = -986.5714
$ws.Range("H132").Value = 23640.826
$ws.Range("I132").Value = 1525.1177
$ws.Range("K132").Value = 4575.3531
$ws.Range("M132").Value = -2045.3531
$ws.Range("H136").Value = 2436.8147
$ws.Range("I136").Value = 1198.5714
$ws.Range("K136").Value = 3595.7142
$ws.Range("M136").Value = -1045.7142

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1999.8
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 2999
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 2999
$ws.Range("M99").Value = -252
$ws.Range("N99").Value = -5995

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12326.818
$ws.Range("J58").Value = 72228.57000000001
$ws.Range("L58").Value = 72228.57000000001
$ws.Range("N58").Value = -72634.57000000001
$ws.Range("H99").Value = 6112.5
$ws.Range("I99").Value = 4328.5713
$ws.Range("J99").Value = 7500
$ws.Range("K99").Value = 4328.5713
$ws.Range("L99").Value = 7500
$ws.Range("M99").Value = -2830.5713
$ws.Range("N99").Value = -10496
$ws.Range("H126").Value = 6112.5
$ws.Range("I126").Value = 4328.5713
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 12985.7139
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -10515.7139
$ws.Range("N126").Value = -27440
$ws.Range("H132").Value = 28552.65
$ws.Range("I132").Value = 46792.09
$ws.Range("K132").Value = 140376.27
$ws.Range("M132").Value = -137846.27
$ws.Range("H134").Value = 947.92
$ws.Range("I134").Value = 935
$ws.Range("K134").Value = 2805
$ws.Range("M134").Value = -270
$ws.Range("H136").Value = 12326.818
$ws.Range("J136").Value = 72228.57000000001
$ws.Range("L136").Value = 216685.71
$ws.Range("N136").Value = -221785.71

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 469.54544
$ws.Range("I26").Value = 153
$ws.Range("J26").Value = 733.3333
$ws.Range("K26").Value = 459
$ws.Range("L26").Value = 2199.9999
$ws.Range("M26").Value = -171
$ws.Range("N26").Value = -2775.9999
$ws.Range("H68").Value = 5072
$ws.Range("I68").Value = 580.7273
$ws.Range("J68").Value = 8365.6
$ws.Range("K68").Value = 1742.1819
$ws.Range("L68").Value = 25096.8
$ws.Range("M68").Value = -931.1819
$ws.Range("N68").Value = -26718.8
$ws.Range("H70").Value = 3352.7058
$ws.Range("I70").Value = 2446.5454
$ws.Range("J70").Value = 5014
$ws.Range("K70").Value = 7339.6362
$ws.Range("L70").Value = 15042
$ws.Range("M70").Value = -7024.6362
$ws.Range("N70").Value = -15672
$ws.Range("H71").Value = 5072
$ws.Range("I71").Value = 580.7273
$ws.Range("J71").Value = 8365.6
$ws.Range("K71").Value = 5226.545700000001
$ws.Range("L71").Value = 75290.40000000001
$ws.Range("M71").Value = -1170.545700000001
$ws.Range("N71").Value = -83402.40000000001
$ws.Range("H73").Value = 3352.7058
$ws.Range("I73").Value = 2446.5454
$ws.Range("J73").Value = 5014
$ws.Range("K73").Value = 7339.6362
$ws.Range("L73").Value = 15042
$ws.Range("M73").Value = -6247.6362
$ws.Range("N73").Value = -17226
$ws.Range("H107").Value = 4366.815
$ws.Range("I107").Value = 11607.333
$ws.Range("J107").Value = 746.55554
$ws.Range("K107").Value = 34821.999
$ws.Range("L107").Value = 2239.66662
$ws.Range("M107").Value = -32901.999
$ws.Range("N107").Value = -6079.66662
$ws.Range("H131").Value = 787.29
$ws.Range("I131").Value = 236.66667
$ws.Range("J131").Value = 804.3196
$ws.Range("K131").Value = 710.00001
$ws.Range("L131").Value = 2412.9588
$ws.Range("M131").Value = 4329.99999
$ws.Range("N131").Value = -12492.9588
$ws.Range("H132").Value = 935.5
$ws.Range("J132").Value = 873.75
$ws.Range("L132").Value = 7863.75
$ws.Range("N132").Value = -12923.75

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1735.5714
$ws.Range("I107").Value = 419.8
$ws.Range("J107").Value = 5025
$ws.Range("K107").Value = 419.8
$ws.Range("L107").Value = 5025
$ws.Range("M107").Value = 1500.2
$ws.Range("N107").Value = -8865
$ws.Range("H132").Value = 20328.467
$ws.Range("I132").Value = 4082.8
$ws.Range("K132").Value = 12248.4
$ws.Range("M132").Value = -9718.400000000001

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4498.963
$ws.Range("I7").Value = 5186.1113
$ws.Range("J7").Value = 3124.6667
$ws.Range("K7").Value = 5186.1113
$ws.Range("L7").Value = 3124.6667
$ws.Range("M7").Value = -5074.1113
$ws.Range("N7").Value = -3348.6667
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -20980
$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 20000
$ws.Range("N52").Value = -20466
$ws.Range("H126").Value = 4498.963
$ws.Range("I126").Value = 5186.1113
$ws.Range("J126").Value = 3124.6667
$ws.Range("K126").Value = 15558.3339
$ws.Range("L126").Value = 9374.000100000001
$ws.Range("M126").Value = -13088.3339
$ws.Range("N126").Value = -14314.0001
$ws.Range("H132").Value = 1703.375
$ws.Range("I132").Value = 1111.9445
$ws.Range("K132").Value = 3335.8335
$ws.Range("M132").Value = -805.8335000000002

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2704057.8
$ws.Range("I113").Value = 1631.125
$ws.Range("K113").Value = 4893.375
$ws.Range("M113").Value = -2723.375
